{"js": "const replacements = [\n  [\"2024-03-11 Monday\", \"2024-03-12 Tuesday\"],\n  [\"172\u00f73=57, 1\", \"912\u00f76=152, 0\"],\n  [\"446\u00f74=111, 2\", \"308\u00f75=61, 3\"],\n  [\"690\u00f77=98, 4\", \"713\u00f77=101, 6\"],\n  [\"966\u00f78=120, 6\", \"424\u00f79=47, 1\"],\n  [\"746\u00f75=149, 1\", \"779\u00f77=111, 2\"],\n  [\"484\u00f79=53, 7\", \"785\u00f74=196, 1\"],\n  [\"992\u00f79=110, 2\", \"428\u00f73=142, 2\"],\n  [\"997\u00f77=142, 3\", \"901\u00f78=112, 5\"],\n  [\"576\u00f72=288, 0\", \"218\u00f79=24, 2\"],\n  [\"688\u00f72=344, 0\", \"208\u00f74=52, 0\"],\n  [\"465\u00f79=51, 6\", \"374\u00f72=187, 0\"],\n  [\"715\u00f77=102, 1\", \"109\u00f79=12, 1\"],\n  [\"774\u00f76=129, 0\", \"604\u00f72=302, 0\"],\n  [\"259\u00f74=64, 3\", \"598\u00f73=199, 1\"],\n  [\"775\u00f73=258, 1\", \"998\u00f72=499, 0\"],\n  [\"333\u00f76=55, 3\", \"889\u00f77=127, 0\"],\n  [\"399\u00f78=49, 7\", \"838\u00f72=419, 0\"],\n  [\"589\u00f73=196, 1\", \"908\u00f78=113, 4\"],\n  [\"402\u00f74=100, 2\", \"856\u00f78=107, 0\"],\n  [\"195\u00f75=39, 0\", \"494\u00f76=82, 2\"],\n  [\"205\u00f74=51, 1\", \"266\u00f75=53, 1\"],\n  [\"414\u00f78=51, 6\", \"422\u00f77=60, 2\"],\n  [\"617\u00f76=102, 5\", \"466\u00f75=93, 1\"],\n  [\"850\u00f76=141, 4\", \"168\u00f75=33, 3\"],\n  [\"573\u00f76=95, 3\", \"569\u00f76=94, 5\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, 'Replace');\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-03-11 Monday\", \"2024-03-12 Tuesday\"),\n    @(\"172\u00f73=57, 1\", \"912\u00f76=152, 0\"),\n    @(\"446\u00f74=111, 2\", \"308\u00f75=61, 3\"),\n    @(\"690\u00f77=98, 4\", \"713\u00f77=101, 6\"),\n    @(\"966\u00f78=120, 6\", \"424\u00f79=47, 1\"),\n    @(\"746\u00f75=149, 1\", \"779\u00f77=111, 2\"),\n    @(\"484\u00f79=53, 7\", \"785\u00f74=196, 1\"),\n    @(\"992\u00f79=110, 2\", \"428\u00f73=142, 2\"),\n    @(\"997\u00f77=142, 3\", \"901\u00f78=112, 5\"),\n    @(\"576\u00f72=288, 0\", \"218\u00f79=24, 2\"),\n    @(\"688\u00f72=344, 0\", \"208\u00f74=52, 0\"),\n    @(\"465\u00f79=51, 6\", \"374\u00f72=187, 0\"),\n    @(\"715\u00f77=102, 1\", \"109\u00f79=12, 1\"),\n    @(\"774\u00f76=129, 0\", \"604\u00f72=302, 0\"),\n    @(\"259\u00f74=64, 3\", \"598\u00f73=199, 1\"),\n    @(\"775\u00f73=258, 1\", \"998\u00f72=499, 0\"),\n    @(\"333\u00f76=55, 3\", \"889\u00f77=127, 0\"),\n    @(\"399\u00f78=49, 7\", \"838\u00f72=419, 0\"),\n    @(\"589\u00f73=196, 1\", \"908\u00f78=113, 4\"),\n    @(\"402\u00f74=100, 2\", \"856\u00f78=107, 0\"),\n    @(\"195\u00f75=39, 0\", \"494\u00f76=82, 2\"),\n    @(\"205\u00f74=51, 1\", \"266\u00f75=53, 1\"),\n    @(\"414\u00f78=51, 6\", \"422\u00f77=60, 2\"),\n    @(\"617\u00f76=102, 5\", \"466\u00f75=93, 1\"),\n    @(\"850\u00f76=141, 4\", \"168\u00f75=33, 3\"),\n    @(\"573\u00f76=95, 3\", \"569\u00f76=94, 5\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # wdFindContinue=1, wdReplaceAll=2\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}"}
